$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 25 new rows at row 163, pushing the existing "*ない形" content
# (previously rows 163-187) down to rows 188-212.
$ws.Rows("163:187").Insert()

# Fill the newly inserted rows with the new content: end of 第19課文型
# continuation, 第20課文型, 第21課文型 placeholders, and full 第22課文型.
$ws.Range("A163").Value = '*第20課文型'

$ws.Range("A164").Value = '待編輯'

$ws.Range("A165").Value = '*第21課文型'

$ws.Range("A166").Value = '待編輯'

$ws.Range("A167").Value = '*第22課文型'

$ws.Range("A168").Value = 'これは ミラーさんが つくった ケーキです'
$ws.Range("B168").Value = 'これは ミラーさんが 作った ケーキです'
$ws.Range("C168").Value = '這是米勒先生做的蛋糕'

$ws.Range("A169").Value = 'あそこに いる ひとは ミラーさんです'
$ws.Range("B169").Value = 'あそこに いる 人は ミラーさんです'
$ws.Range("C169").Value = '那邊那位就是米勒先生'

$ws.Range("A170").Value = 'きのう ならった ことばを わすれました'
$ws.Range("B170").Value = 'きのう 習った ことばを 忘れました'
$ws.Range("C170").Value = '我忘了昨天學的詞了'

$ws.Range("A171").Value = 'かいものに いく じかんが ありません'
$ws.Range("B171").Value = '買い物に 行く 時間が ありません'
$ws.Range("C171").Value = '我沒時間去購物'

$ws.Range("A172").Value = 'これは ばんり の ちょうじょうで とった しゃしんです'
$ws.Range("B172").Value = 'これは 万里 の 長城で 撮った 写真です'
$ws.Range("C172").Value = '這是我在萬里長城拍的照片'

$ws.Range("A173").Value = '…そうですか すごいですね'
$ws.Range("C173").Value = '我明白了,真了不起'

$ws.Range("A174").Value = 'カリナさんが かいた えは どれですか'
$ws.Range("B174").Value = 'カリナさんが かいた 絵は どれですか'
$ws.Range("C174").Value = '卡琳娜畫的是哪幅畫？'

$ws.Range("A175").Value = 'あれです。 あの うみの えです。'
$ws.Range("B175").Value = 'あれです。 あの 海の 絵です。'
$ws.Range("C175").Value = '就是那幅畫，畫的是大海。'

$ws.Range("A176").Value = 'あの き ものを きて いる ひとは だれですか。'
$ws.Range("B176").Value = 'あの 着物を 着て いる 人は だれですか。'
$ws.Range("C176").Value = '那位穿和服的人是誰？'

$ws.Range("A177").Value = 'き むらさんです。'
$ws.Range("B177").Value = '木村さんです。'
$ws.Range("C177").Value = '是木村先生。'

$ws.Range("A178").Value = 'やまださん、おくさんに はじめて あった ところは どこですか。'
$ws.Range("B178").Value = '山田さん、奥さんに 初めて 会った 所は どこですか。'
$ws.Range("C178").Value = '山田先生，您和您太太第一次見面是在哪裡？'

$ws.Range("A179").Value = 'おおさかじょうです。'
$ws.Range("B179").Value = '大阪城です。'
$ws.Range("C179").Value = '是在大阪。'

$ws.Range("A180").Value = 'きむらさんと いった コンサートは どうでしたか。'
$ws.Range("B180").Value = '木村さんと 行った コンサートは どうでしたか。'
$ws.Range("C180").Value = '你和木村先生一起去看的演唱會怎麼樣？'

$ws.Range("A181").Value = 'とても よかったです。'
$ws.Range("C181").Value = '真的很好。'

$ws.Range("A182").Value = 'どう しましたか。'
$ws.Range("C182").Value = '怎麼了？'

$ws.Range("A183").Value = 'きのう かった かさを なくしました。'
$ws.Range("B183").Value = 'きのう 買った 傘を なくしました。'
$ws.Range("C183").Value = '我把昨天買的傘弄丟了。'

$ws.Range("A184").Value = 'どんな うちが ほしいですか。'
$ws.Range("B184").Value = 'どんな うちが 欲しいですか。'
$ws.Range("C184").Value = '你想要什麼樣的房子？'

$ws.Range("A185").Value = 'ひろい にわが ある うちが ほしいです。'
$ws.Range("B185").Value = '広い 庭が ある うちが 欲しいです。'
$ws.Range("C185").Value = '我想要一棟有大花園的房子。'

$ws.Range("A186").Value = 'こんばん のみに いきませんか。'
$ws.Range("B186").Value = '今晩 飲みに 行きませんか。'
$ws.Range("C186").Value = '今晚想出去喝一杯嗎？'

$ws.Range("A187").Value = 'すみません。こんばんは ちょっと ともだちに あう やくそくが あります。'
$ws.Range("B187").Value = 'すみません。今晩は ちょっと 友達に 会う 約束が あります。'
$ws.Range("C187").Value = '抱歉，晚上我約了個朋友。'

# Update the active selection to match the edited document state.
$ws.Activate()
$ws.Range("A166").Select()